$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column C (target stored width 10.7109375 chars; COM ColumnWidth is
# quantized to 1/6-character pixel steps by the engine, so 9.833333333333334
# is the input that lands on the closest achievable stored width).
$ws.Columns.Item(3).ColumnWidth = 9.833333333333334

# Update the three data values on row 1
$ws.Range("A1").Value = 156.63574734127766
$ws.Range("B1").Value = 4.9429305972999993
$ws.Range("C1").Value = 1.3585127201565559
